$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4
$ws.Range("A4").Value = 112388103
$ws.Range("B4").Value = 77636
$ws.Range("C4").Value = "Ovaliderad"
$ws.Range("D4").Value = "NT"
$ws.Range("E4").Value = 6425
$ws.Range("F4").Value = "Garnlav"
$ws.Range("G4").Value = "Alectoria sarmentosa"
$ws.Range("H4").Value = "(Ach.) Ach."
$ws.Range("I4").Value = ""
$ws.Range("J4").Value = ""
$ws.Range("K4").Value = ""
$ws.Range("N4").Value = ""
$ws.Range("P4").Value = "Mörtsjöberget, Ång"
$ws.Range("Q4").Value = 557984
$ws.Range("R4").Value = 7069574
$ws.Range("S4").Value = 25
$ws.Range("T4").Value = "Jämtland"
$ws.Range("U4").Value = "Strömsund"
$ws.Range("V4").Value = "Ångermanland"
$ws.Range("W4").Value = "Fjällsjö"
$ws.Range("Y4").NumberFormat = "@"
$ws.Range("Y4").Value = "2023-09-28"
$ws.Range("Y4").Style = "Normal"
$ws.Range("AA4").NumberFormat = "@"
$ws.Range("AA4").Value = "2023-09-28"
$ws.Range("AA4").Style = "Normal"
$ws.Range("AD4").Value = $false
$ws.Range("AE4").Value = $false
$ws.Range("AF4").Value = ""
$ws.Range("AG4").Value = $false
$ws.Range("AT4").Value = ""
$ws.Range("AW4").Value = "Maria Johansson"
$ws.Range("AX4").Value = "Maria Johansson"
$ws.Range("AY4").Value = ""

# Row 5
$ws.Range("A5").Value = 112388110
$ws.Range("B5").Value = 77636
$ws.Range("C5").Value = "Ovaliderad"
$ws.Range("D5").Value = "NT"
$ws.Range("E5").Value = 6425
$ws.Range("F5").Value = "Garnlav"
$ws.Range("G5").Value = "Alectoria sarmentosa"
$ws.Range("H5").Value = "(Ach.) Ach."
$ws.Range("I5").Value = ""
$ws.Range("J5").Value = ""
$ws.Range("K5").Value = ""
$ws.Range("N5").Value = ""
$ws.Range("P5").Value = "Mörtsjöberget, Ång"
$ws.Range("Q5").Value = 557867
$ws.Range("R5").Value = 7069706
$ws.Range("S5").Value = 25
$ws.Range("T5").Value = "Jämtland"
$ws.Range("U5").Value = "Strömsund"
$ws.Range("V5").Value = "Ångermanland"
$ws.Range("W5").Value = "Fjällsjö"
$ws.Range("Y5").NumberFormat = "@"
$ws.Range("Y5").Value = "2023-09-28"
$ws.Range("Y5").Style = "Normal"
$ws.Range("AA5").NumberFormat = "@"
$ws.Range("AA5").Value = "2023-09-28"
$ws.Range("AA5").Style = "Normal"
$ws.Range("AD5").Value = $false
$ws.Range("AE5").Value = $false
$ws.Range("AF5").Value = ""
$ws.Range("AG5").Value = $false
$ws.Range("AT5").Value = ""
$ws.Range("AW5").Value = "Maria Johansson"
$ws.Range("AX5").Value = "Maria Johansson"
$ws.Range("AY5").Value = ""

# Row 6
$ws.Range("A6").Value = 112388101
$ws.Range("B6").Value = 96720
$ws.Range("C6").Value = "Ovaliderad"
$ws.Range("D6").Value = "VU"
$ws.Range("E6").Value = 220787
$ws.Range("F6").Value = "Knärot"
$ws.Range("G6").Value = "Goodyera repens"
$ws.Range("H6").Value = "(L.) R. Br."
$ws.Range("I6").Value = ""
$ws.Range("J6").Value = ""
$ws.Range("K6").Value = ""
$ws.Range("L6").Value = ""
$ws.Range("N6").Value = ""
$ws.Range("P6").Value = "Mörtsjöberget, Ång"
$ws.Range("Q6").Value = 557984
$ws.Range("R6").Value = 7069575
$ws.Range("S6").Value = 25
$ws.Range("T6").Value = "Jämtland"
$ws.Range("U6").Value = "Strömsund"
$ws.Range("V6").Value = "Ångermanland"
$ws.Range("W6").Value = "Fjällsjö"
$ws.Range("Y6").NumberFormat = "@"
$ws.Range("Y6").Value = "2023-09-28"
$ws.Range("Y6").Style = "Normal"
$ws.Range("AA6").NumberFormat = "@"
$ws.Range("AA6").Value = "2023-09-28"
$ws.Range("AA6").Style = "Normal"
$ws.Range("AD6").Value = $false
$ws.Range("AE6").Value = $false
$ws.Range("AF6").Value = ""
$ws.Range("AG6").Value = $false
$ws.Range("AT6").Value = ""
$ws.Range("AW6").Value = "Maria Johansson"
$ws.Range("AX6").Value = "Maria Johansson"
$ws.Range("AY6").Value = ""

# Row 7
$ws.Range("A7").Value = 112388107
$ws.Range("B7").Value = 81371
$ws.Range("C7").Value = "Ovaliderad"
$ws.Range("D7").Value = "NT"
$ws.Range("E7").Value = 1312
$ws.Range("F7").Value = "Gammelgransskål"
$ws.Range("G7").Value = "Pseudographis pinicola"
$ws.Range("H7").Value = "(Nyl.) Rehm"
$ws.Range("I7").Value = ""
$ws.Range("J7").Value = ""
$ws.Range("K7").Value = ""
$ws.Range("N7").Value = ""
$ws.Range("P7").Value = "Mörtsjöberget, Ång"
$ws.Range("Q7").Value = 557867
$ws.Range("R7").Value = 7069709
$ws.Range("S7").Value = 25
$ws.Range("T7").Value = "Jämtland"
$ws.Range("U7").Value = "Strömsund"
$ws.Range("V7").Value = "Ångermanland"
$ws.Range("W7").Value = "Fjällsjö"
$ws.Range("Y7").NumberFormat = "@"
$ws.Range("Y7").Value = "2023-09-28"
$ws.Range("Y7").Style = "Normal"
$ws.Range("AA7").NumberFormat = "@"
$ws.Range("AA7").Value = "2023-09-28"
$ws.Range("AA7").Style = "Normal"
$ws.Range("AD7").Value = $false
$ws.Range("AE7").Value = $false
$ws.Range("AF7").Value = ""
$ws.Range("AG7").Value = $false
$ws.Range("AT7").Value = ""
$ws.Range("AW7").Value = "Maria Johansson"
$ws.Range("AX7").Value = "Maria Johansson"
$ws.Range("AY7").Value = ""

# Row 8
$ws.Range("A8").Value = 112388115
$ws.Range("B8").Value = 96720
$ws.Range("C8").Value = "Ovaliderad"
$ws.Range("D8").Value = "VU"
$ws.Range("E8").Value = 220787
$ws.Range("F8").Value = "Knärot"
$ws.Range("G8").Value = "Goodyera repens"
$ws.Range("H8").Value = "(L.) R. Br."
$ws.Range("I8").Value = ""
$ws.Range("J8").Value = ""
$ws.Range("K8").Value = ""
$ws.Range("L8").Value = ""
$ws.Range("N8").Value = ""
$ws.Range("P8").Value = "Mörtsjöberget, Ång"
$ws.Range("Q8").Value = 557811
$ws.Range("R8").Value = 7069647
$ws.Range("S8").Value = 25
$ws.Range("T8").Value = "Jämtland"
$ws.Range("U8").Value = "Strömsund"
$ws.Range("V8").Value = "Ångermanland"
$ws.Range("W8").Value = "Fjällsjö"
$ws.Range("Y8").NumberFormat = "@"
$ws.Range("Y8").Value = "2023-09-28"
$ws.Range("Y8").Style = "Normal"
$ws.Range("AA8").NumberFormat = "@"
$ws.Range("AA8").Value = "2023-09-28"
$ws.Range("AA8").Style = "Normal"
$ws.Range("AD8").Value = $false
$ws.Range("AE8").Value = $false
$ws.Range("AF8").Value = ""
$ws.Range("AG8").Value = $false
$ws.Range("AT8").Value = ""
$ws.Range("AW8").Value = "Maria Johansson"
$ws.Range("AX8").Value = "Maria Johansson"
$ws.Range("AY8").Value = ""

# Row 9
$ws.Range("A9").Value = 112388117
$ws.Range("B9").Value = 77636
$ws.Range("C9").Value = "Ovaliderad"
$ws.Range("D9").Value = "NT"
$ws.Range("E9").Value = 6425
$ws.Range("F9").Value = "Garnlav"
$ws.Range("G9").Value = "Alectoria sarmentosa"
$ws.Range("H9").Value = "(Ach.) Ach."
$ws.Range("I9").Value = ""
$ws.Range("J9").Value = ""
$ws.Range("K9").Value = ""
$ws.Range("N9").Value = ""
$ws.Range("P9").Value = "Mörtsjöberget, Ång"
$ws.Range("Q9").Value = 557810
$ws.Range("R9").Value = 7069645
$ws.Range("S9").Value = 25
$ws.Range("T9").Value = "Jämtland"
$ws.Range("U9").Value = "Strömsund"
$ws.Range("V9").Value = "Ångermanland"
$ws.Range("W9").Value = "Fjällsjö"
$ws.Range("Y9").NumberFormat = "@"
$ws.Range("Y9").Value = "2023-09-28"
$ws.Range("Y9").Style = "Normal"
$ws.Range("AA9").NumberFormat = "@"
$ws.Range("AA9").Value = "2023-09-28"
$ws.Range("AA9").Style = "Normal"
$ws.Range("AD9").Value = $false
$ws.Range("AE9").Value = $false
$ws.Range("AF9").Value = ""
$ws.Range("AG9").Value = $false
$ws.Range("AT9").Value = ""
$ws.Range("AW9").Value = "Maria Johansson"
$ws.Range("AX9").Value = "Maria Johansson"
$ws.Range("AY9").Value = ""

